# Add a new "Correction " column (N) to the Card23 sheet, right after the
# "Event" column (M), mirroring the layout already used on the sibling
# "CardNN" sheets (e.g. Card24 / Card22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card23")

# 1) Header: "Event " (trailing space) becomes "Event" (no trailing space).
$ws.Range("M1").Value = "Event"

# 2) Copy the header's formatting (bold font, border, centered alignment)
#    onto the new column's header cell, then give it its text.
$ws.Range("M1").Copy()
$ws.Range("N1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N1").Value = "Correction "

# 3) Data rows 2-12: the existing (empty) "Event" cells become the literal
#    text "nan", matching every other text column on this sheet, and the
#    new "Correction" column is added alongside them (left blank).
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Value = "nan"   # column M
    $ws.Cells.Item($r, 14).Value = ""      # column N (new, blank)
}
